$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.101.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.658.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.98%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.657.06"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.13%  "
$ws.Range("E11").Value = "  +3.52%  "
$ws.Range("E12").Value = "  +1.93%  "
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.101.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.049.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.79%  "
$ws.Range("E17").Value = "  +1.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.654.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "350.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.61%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.87%  "
$ws.Range("E22").Value = "  +4.03%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.57%  "
$ws.Range("E25").Value = "  +2.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.764.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.12%  "
$ws.Range("E27").Value = "  +2.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.995"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0807"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  +7.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.80%  "
$ws.Range("E34").Value = "  +3.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.966"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.24%  "
$ws.Range("E37").Value = "  +3.68%  "
$ws.Range("E38").Value = "  +2.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("E40").Value = "  +3.07%  "
$ws.Range("E41").Value = "  +5.36%  "
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "277.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0985"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("E46").Value = "  +1.71%  "
$ws.Range("E47").Value = "  +5.68%  "
$ws.Range("E48").Value = "  -2.02%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0230"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.51%  "
